# Insert a new weekly record at row 57 for "Arveja Verde" (Macroferia Regional de Talca).
# This shifts the existing rows 57-165 down to 58-166 (Excel preserves all the
# other cell values/formatting automatically when inserting a whole row), and
# we then populate the freshly inserted row 57 with the new data point.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 57, pushing rows 57:165 to 58:166.
$ws.Rows.Item(57).Insert()

# Populate the new row 57 with the new data record.
$ws.Range("A57").Value = 5
$ws.Range("B57").Value = "Macroferia Regional de Talca"
$ws.Range("C57").Value = "Maule"
$ws.Range("D57").Value = 45259
$ws.Range("E57").Value = 7
$ws.Range("F57").Value = 100112022
$ws.Range("G57").Value = "Arveja Verde"
$ws.Range("H57").Value = "Sin especificar"
$ws.Range("I57").Value = "Primera"
$ws.Range("J57").Value = 400
$ws.Range("K57").Value = 15000
$ws.Range("L57").Value = 18000
$ws.Range("M57").Value = 16500
$ws.Range("N57").Value = "`$/saco 25 kilos"
$ws.Range("O57").Value = "Región del Maule"
$ws.Range("P57").Value = 660
$ws.Range("Q57").Value = 25
$ws.Range("R57").Value = "Hortaliza"

# Keep the date column formatted the same way as the rest of column D.
$ws.Range("D57").NumberFormat = $ws.Range("D58").NumberFormat
